# aggiornamento 15, 16, 17 marzo
# Append three new daily rows (227-229) to the single data sheet, continuing
# the existing date/number sequence, and copy the date-cell formatting from
# the last existing row (A226) onto the new date cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data: date serial, nuovi pos., somma mobile 7gg., somma mobile 7gg. per 100mila abitanti
$newRows = @(
    @(44301, 1, 33, 218.7313581228873),
    @(44302, 2, 32, 212.1031351494664),
    @(44303, 7, 28, 185.5902432557831)
)

$startRow = 227

for ($i = 0; $i -lt $newRows.Length; $i++) {
    $r = $startRow + $i
    $data = $newRows[$i]

    # Match the formatting used by the preceding date cell (A column).
    $ws.Cells.Item($r - 1, 1).Copy()
    $ws.Cells.Item($r, 1).PasteSpecial(-4122)  # xlPasteFormats

    $ws.Cells.Item($r, 1).Value = $data[0]
    $ws.Cells.Item($r, 2).Value = $data[1]
    $ws.Cells.Item($r, 3).Value = $data[2]
    $ws.Cells.Item($r, 4).Value = $data[3]
}
